$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at row 78. This pushes the existing rows
# 78..201 down to 79..202 (preserving all their values/styles) and
# grows the sheet's used range from A1:R201 to A1:R202, matching the
# target diff.
$ws.Rows(78).Insert()

# Populate the newly-inserted row 78 with the new weekly data point.
$ws.Cells.Item(78, 1).Value = 4
$ws.Cells.Item(78, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(78, 3).Value = "Los Lagos"
$ws.Cells.Item(78, 4).Value = 44495
$ws.Cells.Item(78, 5).Value = 10
$ws.Cells.Item(78, 6).Value = 100114014
$ws.Cells.Item(78, 7).Value = "Betarraga"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 1200
$ws.Cells.Item(78, 11).Value = 1000
$ws.Cells.Item(78, 12).Value = 1200
$ws.Cells.Item(78, 13).Value = 1100
$ws.Cells.Item(78, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(78, 15).Value = "Región del Maule"
$ws.Cells.Item(78, 16).Value = 220
$ws.Cells.Item(78, 17).Value = 5
$ws.Cells.Item(78, 18).Value = "Hortaliza"
